# This script applies the odds updates captured in the commit diff for
# "Jogos_da_Semana_FlashScore_2024-12-12.xlsx" (a weekly football odds export).
# Each football match occupies one row (rows 2-11); every Range().Value
# assignment below updates a single odds cell to its new value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("O2").Value = 1.17
$ws.Range("P2").Value = 5
$ws.Range("Q2").Value = 1.57
$ws.Range("R2").Value = 2.35
# Row 4
$ws.Range("AC4").Value = 12.5
$ws.Range("AD4").Value = 6.6
$ws.Range("AH4").Value = 11.25
$ws.Range("AT4").Value = 3
$ws.Range("AX4").Value = 5.9
$ws.Range("AY4").Value = 22
$ws.Range("H4").Value = 3.75
$ws.Range("I4").Value = 3.95
$ws.Range("J4").Value = 2.22
$ws.Range("M4").Value = 1.02
$ws.Range("N4").Value = 9.8
$ws.Range("O4").Value = 1.19
$ws.Range("P4").Value = 4.21
$ws.Range("Q4").Value = 1.65
$ws.Range("R4").Value = 2
$ws.Range("U4").Value = 1.63
$ws.Range("V4").Value = 2.15
$ws.Range("W4").Value = 7
$ws.Range("X4").Value = 7.4
$ws.Range("Z4").Value = 11
# Row 5
$ws.Range("O5").Value = 1.18
$ws.Range("P5").Value = 4.5
$ws.Range("Q5").Value = 1.65
$ws.Range("R5").Value = 2.2
# Row 6
$ws.Range("AD6").Value = 7.2
$ws.Range("AE6").Value = 15.5
$ws.Range("AF6").Value = 70
$ws.Range("AG6").Value = 500
$ws.Range("AH6").Value = 14.5
$ws.Range("AI6").Value = 32
$ws.Range("AJ6").Value = 16
$ws.Range("AK6").Value = 100
$ws.Range("AN6").Value = 3.4
$ws.Range("AO6").Value = 8
$ws.Range("AP6").Value = 18
$ws.Range("AQ6").Value = 27
$ws.Range("AR6").Value = 60
$ws.Range("AT6").Value = 2.52
$ws.Range("AU6").Value = 7.5
$ws.Range("AX6").Value = 6.6
$ws.Range("AY6").Value = 29
$ws.Range("G6").Value = 1.62
$ws.Range("I6").Value = 5.1
$ws.Range("K6").Value = 2.12
$ws.Range("L6").Value = 5.1
$ws.Range("N6").Value = 10.5
$ws.Range("S6").Value = 1.39
$ws.Range("T6").Value = 2.57
$ws.Range("W6").Value = 6.7
$ws.Range("X6").Value = 7.5
$ws.Range("Z6").Value = 12
# Row 7
$ws.Range("AA7").Value = 15
$ws.Range("AC7").Value = 8.25
$ws.Range("AD7").Value = 6.7
$ws.Range("AJ7").Value = 17
$ws.Range("AK7").Value = 110
$ws.Range("AL7").Value = 60
$ws.Range("AO7").Value = 7.9
$ws.Range("AP7").Value = 18
$ws.Range("AQ7").Value = 27
$ws.Range("AT7").Value = 2.52
$ws.Range("AU7").Value = 7.4
$ws.Range("AY7").Value = 30
$ws.Range("BC7").Value = 400
$ws.Range("H7").Value = 3.35
$ws.Range("I7").Value = 5.3
$ws.Range("J7").Value = 2.2
$ws.Range("K7").Value = 2.12
$ws.Range("M7").Value = 1.01
$ws.Range("N7").Value = 8
$ws.Range("O7").Value = 1.35
$ws.Range("P7").Value = 2.72
$ws.Range("Q7").Value = 2.02
$ws.Range("R7").Value = 1.62
$ws.Range("T7").Value = 2.57
$ws.Range("W7").Value = 5.6
$ws.Range("X7").Value = 6.9
$ws.Range("Z7").Value = 12.5
# Row 9
$ws.Range("M9").Value = 1.06
$ws.Range("N9").Value = 10
$ws.Range("O9").Value = 1.33
$ws.Range("P9").Value = 3.25
# Row 10
$ws.Range("AD10").Value = 7
$ws.Range("AE10").Value = 13
$ws.Range("AG10").Value = 151
$ws.Range("AJ10").Value = 13
$ws.Range("AL10").Value = 29
$ws.Range("AN10").Value = 4
$ws.Range("AO10").Value = 9.5
$ws.Range("AS10").Value = 101
$ws.Range("AU10").Value = 7.5
$ws.Range("AX10").Value = 6
$ws.Range("AY10").Value = 21
$ws.Range("AZ10").Value = 26
$ws.Range("BA10").Value = 67
$ws.Range("G10").Value = 1.8
$ws.Range("H10").Value = 3.6
$ws.Range("I10").Value = 4
$ws.Range("J10").Value = 2.38
$ws.Range("L10").Value = 4.33
$ws.Range("U10").Value = 1.67
$ws.Range("V10").Value = 2.1
$ws.Range("W10").Value = 8.5
$ws.Range("X10").Value = 9.5
# Row 11
$ws.Range("AC11").Value = 6.8
$ws.Range("AD11").Value = 7.1
$ws.Range("AI11").Value = 6.9
$ws.Range("AK11").Value = 11.75
$ws.Range("AN11").Value = 6.6
$ws.Range("AO11").Value = 30
$ws.Range("AT11").Value = 2.67
$ws.Range("AU11").Value = 8
$ws.Range("AY11").Value = 7.9
$ws.Range("BA11").Value = 27
$ws.Range("BB11").Value = 65
$ws.Range("G11").Value = 4.9
$ws.Range("I11").Value = 1.62
$ws.Range("K11").Value = 2.15
$ws.Range("L11").Value = 2.2
$ws.Range("N11").Value = 6.8
$ws.Range("O11").Value = 1.35
$ws.Range("P11").Value = 2.92
$ws.Range("S11").Value = 1.42
$ws.Range("T11").Value = 2.67
$ws.Range("W11").Value = 11.75
$ws.Range("X11").Value = 28
